# Refactor open account test using page objects model (#2)
#
# - "OpenAccountTest" sheet: customer name test data updated from
#   "Joao Silva" to "Harry Potter" (cell A2).
# - The workbook's active/selected sheet moves from "AddCustomerTest"
#   (tab 0) to "OpenAccountTest" (tab 1), and the live selection on
#   "OpenAccountTest" moves from E16 to F9.

$wb = $excel.ActiveWorkbook

$openAccountSheet = $wb.Worksheets.Item("OpenAccountTest")

# Update the test data: customer name used by the Open Account test.
$openAccountSheet.Range("A2").Value = "Harry Potter"

# Make "OpenAccountTest" the active sheet (was "AddCustomerTest"),
# and move its selection/active cell to F9 (was E16).
$openAccountSheet.Activate() | Out-Null
$openAccountSheet.Range("F9").Select() | Out-Null
